$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: E2 gains the course text (was "-")
$ws.Range("E2").Value = "MCT-2A-Eletrônica analóg. e de potência"

# Row 4: D4 loses the course text, becomes "-"
$ws.Range("D4").Value = "-"

# Row 6: D6 loses the course text (becomes "-"), E6 gains it
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "MCT-2A-Eletrônica analóg. e de potência"

# Row 20: B20 loses the course text (becomes "-"), E20 gains it
$ws.Range("B20").Value = "-"
$ws.Range("E20").Value = "ELM-2NA-Circuitos Elétricos 2"

# Row 21: B21 loses the course text (becomes "-"), E21 gains it
$ws.Range("B21").Value = "-"
$ws.Range("E21").Value = "ELM-2NA-Circuitos Elétricos 2"
